$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 4257.9287
$ws.Range("I28").Value = 4200.923
$ws.Range("J28").Value = 4999
$ws.Range("K28").Value = 4200.923
$ws.Range("L28").Value = 4999
$ws.Range("M28").Value = -3715.923
$ws.Range("N28").Value = -5969

$ws.Range("H61").Value = 5400.8887
$ws.Range("I61").Value = 5400.8887
$ws.Range("K61").Value = 16202.6661
$ws.Range("M61").Value = -16030.6661

$ws.Range("H80").Value = 1208.45
$ws.Range("J80").Value = 1582.2
$ws.Range("L80").Value = 4746.6
$ws.Range("N80").Value = -6742.6

$ws.Range("H82").Value = 260.66666
$ws.Range("I82").Value = 260.66666
$ws.Range("K82").Value = 781.9999799999999
$ws.Range("M82").Value = -375.9999799999999

$ws.Range("H83").Value = 1208.45
$ws.Range("J83").Value = 1582.2
$ws.Range("L83").Value = 14239.8
$ws.Range("N83").Value = -24223.8

$ws.Range("H85").Value = 260.66666
$ws.Range("I85").Value = 260.66666
$ws.Range("K85").Value = 781.9999799999999
$ws.Range("M85").Value = 622.0000200000001

$ws.Range("H86").Value = 5701.8335
$ws.Range("I86").Value = 5134.6665
$ws.Range("J86").Value = 6269
$ws.Range("K86").Value = 5134.6665
$ws.Range("L86").Value = 6269
$ws.Range("M86").Value = -4011.6665
$ws.Range("N86").Value = -8515

$ws.Range("H88").Value = 3133
$ws.Range("J88").Value = 3216
$ws.Range("L88").Value = 3216
$ws.Range("N88").Value = -4028

$ws.Range("H89").Value = 5701.8335
$ws.Range("I89").Value = 5134.6665
$ws.Range("J89").Value = 6269
$ws.Range("K89").Value = 25673.3325
$ws.Range("L89").Value = 31345
$ws.Range("M89").Value = -20057.3325
$ws.Range("N89").Value = -42577

$ws.Range("H91").Value = 3133
$ws.Range("J91").Value = 3216
$ws.Range("L91").Value = 3216
$ws.Range("N91").Value = -6024

$ws.Range("H132").Value = 4741.161
$ws.Range("I132").Value = 4112.9585
$ws.Range("J132").Value = 6895
$ws.Range("K132").Value = 12338.8755
$ws.Range("L132").Value = 20685
$ws.Range("M132").Value = -9808.875499999998
$ws.Range("N132").Value = -25745

$ws.Range("H138").Value = 2644.326
$ws.Range("I138").Value = 1366.0435
$ws.Range("J138").Value = 3089.7878
$ws.Range("K138").Value = 4098.1305
$ws.Range("L138").Value = 9269.3634
$ws.Range("M138").Value = 1041.8695
$ws.Range("N138").Value = -19549.3634

$ws.Range("H141").Value = 3173.3572
$ws.Range("I141").Value = 3175.5386
$ws.Range("K141").Value = 9526.6158
$ws.Range("M141").Value = -4346.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3579
$ws.Range("I88").Value = 3088
$ws.Range("K88").Value = 3088
$ws.Range("M88").Value = -2682

$ws.Range("H91").Value = 3579
$ws.Range("I91").Value = 3088
$ws.Range("K91").Value = 3088
$ws.Range("M91").Value = -1684

$ws.Range("H122").Value = 4597
$ws.Range("I122").Value = 2516.4
$ws.Range("K122").Value = 7549.200000000001
$ws.Range("M122").Value = -5099.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2003.5
$ws.Range("I105").Value = 1955.7587
$ws.Range("K105").Value = 1955.7587
$ws.Range("M105").Value = -208.7587000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2706.3333
$ws.Range("I16").Value = 2738.1428
$ws.Range("K16").Value = 2738.1428
$ws.Range("M16").Value = -2451.1428

$ws.Range("H31").Value = 419411.62
$ws.Range("I31").Value = 1002200.9
$ws.Range("J31").Value = 3133.5715
$ws.Range("K31").Value = 1002200.9
$ws.Range("L31").Value = 3133.5715
$ws.Range("M31").Value = -1001905.9
$ws.Range("N31").Value = -3723.5715

$ws.Range("H34").Value = 419411.62
$ws.Range("I34").Value = 1002200.9
$ws.Range("J34").Value = 3133.5715
$ws.Range("K34").Value = 1002200.9
$ws.Range("L34").Value = 3133.5715
$ws.Range("M34").Value = -1001998.9
$ws.Range("N34").Value = -3537.5715

$ws.Range("H62").Value = 1724.2222
$ws.Range("I62").Value = 1901.8
$ws.Range("J62").Value = 1502.25
$ws.Range("K62").Value = 1901.8
$ws.Range("L62").Value = 1502.25
$ws.Range("M62").Value = -1277.8
$ws.Range("N62").Value = -2750.25

$ws.Range("H65").Value = 1724.2222
$ws.Range("I65").Value = 1901.8
$ws.Range("J65").Value = 1502.25
$ws.Range("K65").Value = 9509
$ws.Range("L65").Value = 7511.25
$ws.Range("M65").Value = -6389
$ws.Range("N65").Value = -13751.25

$ws.Range("H68").Value = 37056.1
$ws.Range("J68").Value = 37056.1
$ws.Range("L68").Value = 37056.1
$ws.Range("N68").Value = -38554.1

$ws.Range("H71").Value = 37056.1
$ws.Range("J71").Value = 37056.1
$ws.Range("L71").Value = 111168.3
$ws.Range("N71").Value = -118656.3

$ws.Range("H113").Value = 2706.3333
$ws.Range("I113").Value = 2738.1428
$ws.Range("K113").Value = 2738.1428
$ws.Range("M113").Value = -568.1428000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 525.75
$ws.Range("I108").Value = 525.75
$ws.Range("K108").Value = 1577.25
$ws.Range("M108").Value = 1302.75

$ws.Range("H121").Value = 76707.36
$ws.Range("I121").Value = 5440.6665
$ws.Range("J121").Value = 96143.73
$ws.Range("K121").Value = 16321.9995
$ws.Range("L121").Value = 288431.19
$ws.Range("M121").Value = -15011.9995
$ws.Range("N121").Value = -291051.19

$ws.Range("H131").Value = 7693703
$ws.Range("J131").Value = 1446.65
$ws.Range("L131").Value = 4339.950000000001
$ws.Range("N131").Value = -14419.95

$ws.Range("H132").Value = 4691.184
$ws.Range("I132").Value = 5170.1934
$ws.Range("K132").Value = 46531.7406
$ws.Range("M132").Value = -44001.7406

$ws.Range("H139").Value = 2370.625
$ws.Range("I139").Value = 1625.7222
$ws.Range("K139").Value = 4877.1666
$ws.Range("M139").Value = 262.8334000000004

$ws.Range("H141").Value = 3353
$ws.Range("I141").Value = 3353
$ws.Range("K141").Value = 10059
$ws.Range("M141").Value = -4879

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4009.3845
$ws.Range("I80").Value = 3198
$ws.Range("J80").Value = 4252.8
$ws.Range("K80").Value = 3198
$ws.Range("L80").Value = 4252.8
$ws.Range("M80").Value = -2200
$ws.Range("N80").Value = -6248.8

$ws.Range("H83").Value = 4009.3845
$ws.Range("I83").Value = 3198
$ws.Range("J83").Value = 4252.8
$ws.Range("K83").Value = 15990
$ws.Range("L83").Value = 21264
$ws.Range("M83").Value = -10998
$ws.Range("N83").Value = -31248

$ws.Range("H113").Value = 5497.0967
$ws.Range("I113").Value = 6348.0586
$ws.Range("K113").Value = 6348.0586
$ws.Range("M113").Value = -4178.0586

$ws.Range("H132").Value = 503906
$ws.Range("I132").Value = 503906
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1511718
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4703.615
$ws.Range("I7").Value = 5308.409
$ws.Range("K7").Value = 5308.409
$ws.Range("M7").Value = -5196.409

$ws.Range("H82").Value = 1395.4286
$ws.Range("I82").Value = 1097.6666
$ws.Range("J82").Value = 1792.4445
$ws.Range("K82").Value = 1097.6666
$ws.Range("L82").Value = 1792.4445
$ws.Range("M82").Value = -736.6666
$ws.Range("N82").Value = -2514.4445

$ws.Range("H85").Value = 1395.4286
$ws.Range("I85").Value = 1097.6666
$ws.Range("J85").Value = 1792.4445
$ws.Range("K85").Value = 1097.6666
$ws.Range("L85").Value = 1792.4445
$ws.Range("M85").Value = 150.3334
$ws.Range("N85").Value = -4288.4445

$ws.Range("H126").Value = 4703.615
$ws.Range("I126").Value = 5308.409
$ws.Range("K126").Value = 15925.227
$ws.Range("M126").Value = -13455.227

$ws.Range("H132").Value = 10075.5
$ws.Range("I132").Value = 11079.75
$ws.Range("K132").Value = 33239.25
$ws.Range("M132").Value = -30709.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6196.1333
$ws.Range("I62").Value = 2495.875
$ws.Range("J62").Value = 10425
$ws.Range("K62").Value = 2495.875
$ws.Range("L62").Value = 10425
$ws.Range("M62").Value = -1871.875
$ws.Range("N62").Value = -11673

$ws.Range("H65").Value = 6196.1333
$ws.Range("I65").Value = 2495.875
$ws.Range("J65").Value = 10425
$ws.Range("K65").Value = 12479.375
$ws.Range("L65").Value = 52125
$ws.Range("M65").Value = -9359.375
$ws.Range("N65").Value = -58365

$ws.Range("H81").Value = 3744
$ws.Range("I81").Value = 3744
$ws.Range("K81").Value = 7488
$ws.Range("M81").Value = -6427

$ws.Range("H84").Value = 3744
$ws.Range("I84").Value = 3744
$ws.Range("K84").Value = 37440
$ws.Range("M84").Value = -32136

$ws.Range("H107").Value = 680.875
$ws.Range("I107").Value = 943
$ws.Range("K107").Value = 2829
$ws.Range("M107").Value = -909
